$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the B:G values of each row down into the next row (old row r -> row r+1),
# i.e. row2->row3, row3->row4, ... row10->row11.
# Process from the bottom (r=10) up to the top (r=2) so that a row's original
# value is copied out to the row below before it gets overwritten by the row above it.
for ($r = 10; $r -ge 2; $r--) {
    $nextRow = $r + 1
    $ws.Range("B$nextRow").Value = $ws.Range("B$r").Value2
    $ws.Range("C$nextRow").Value = $ws.Range("C$r").Value2
    $ws.Range("D$nextRow").Value = $ws.Range("D$r").Value2
    $ws.Range("E$nextRow").Value = $ws.Range("E$r").Value2
    $ws.Range("F$nextRow").Value = $ws.Range("F$r").Value2
    $ws.Range("G$nextRow").Value = $ws.Range("G$r").Value2
}

# Row 2 now holds a brand-new leading data point.
$ws.Range("B2").Value = -0.02314597604078636
$ws.Range("C2").Value = 0.3579920056255013
$ws.Range("D2").Value = 0.1782699060034266
$ws.Range("E2").Value = 0.4222202103209018
$ws.Range("F2").Value = 0.4363822494547141
$ws.Range("G2").Value = 15
